# edit.ps1 — reproduce the "Se sube version presentada al profe" commit.
#
# Substance of the change:
#   * Sheet "ROBOS" (2nd sheet): the 15 data rows (rows 2-16, columns A:C)
#     get duplicated again right below the existing data, as rows 17-31
#     (same Location/lat/long values, same relative formatting), and the
#     sheet's selection moves down to E18 (from E17) to follow the new data.
#   * Sheet "CAI" (1st sheet) and the untouched columns of "ROBOS" keep all
#     of their values; only some now-meaningless "apply fill" style flags
#     were cleaned up by Excel on save (no visible formatting change since
#     there never was a fill color — number formats are preserved).
#
# We reproduce the data/content change with the Excel object model.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "CAI"
$ws2 = $wb.Worksheets.Item(2)   # "ROBOS"

# --- ROBOS: duplicate the 15 data rows (A2:C16) into rows 17-31 -----------
$src = $ws2.Range("A2:C16")
$dst = $ws2.Range("A17")
$src.Copy($dst)

# Clear the clipboard marquee / copy mode.
$excel.CutCopyMode = 0

# --- Update the active sheet's selection to follow the new data ----------
$ws2.Activate()
$ws2.Range("E18").Select()

# ROBOS stays the tab that is active/selected when the workbook is saved.
$ws2.Activate()
